$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.019.97"
$ws.Range("E2").Value = "  +0.38%  "
$ws.Range("D3").Value = "1.861.46"
$ws.Range("E3").Value = "  -0.34%  "
$ws.Range("E4").Value = "  +0.06%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "311.62"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -0.21%  "
$ws.Range("E6").Value = "  -0.01%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.5109"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +2.95%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.3814"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +0.24%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.08278"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -6.87%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "1.110"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -0.35%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "41.55"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -0.01%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "6.226"
$c.Style = "Normal"
$ws.Range("E12").Value = "  -1.09%  "
$ws.Range("D13").Value = "1.873.28"
$ws.Range("E13").Value = "  +0.21%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "20.45"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -0.64%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "7.196"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -0.02%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "1.003"
$c.Style = "Normal"
$ws.Range("E17").Value = "  +0.01%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "90.50"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -0.13%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "0.06618"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -0.16%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "17.76"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -0.40%  "
$ws.Range("E21").Value = "  -0.05%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "6.017"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -1.15%  "
$ws.Range("D23").Value = "28.040.16"
$ws.Range("E23").Value = "  +0.38%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "11.07"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -2.29%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "2.227"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -2.44%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "2.567"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +2.67%  "
$ws.Range("D27").Value = "2.080.50"
$ws.Range("E27").Value = "  -0.36%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "157.52"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -0.13%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "20.41"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -1.46%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "124.78"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -1.38%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "0.1061"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +0.80%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "1.037"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -1.19%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "5.608"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +0.78%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "3.602"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +0.40%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "9.616"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +3.90%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "0.06545"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +0.43%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "0.02421"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +1.28%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.2175"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -0.08%  "
$ws.Range("E39").Value = "  +1.11%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "1.243"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -1.55%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.6423"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +1.21%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "11.25"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -3.44%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "4.883"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +0.32%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "0.6083"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +1.92%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "13.07"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -0.60%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "1.277"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -0.74%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "3.650"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -0.66%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "1.982"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +1.18%  "
$ws.Range("E49").Value = "  -0.26%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "120.73"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +0.22%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "79.23"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +1.56%  "
